# Fill in colony counts for Day 19 (column Z), rows 3-34.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    3  = 37
    4  = 37
    5  = 44
    6  = 34
    7  = 106
    8  = 34
    9  = 48
    10 = 312
    11 = 59
    12 = 35
    13 = 72
    14 = 69
    15 = 6
    16 = 32
    17 = 44
    18 = 27
    19 = 79
    20 = 50
    21 = 51
    22 = 47
    23 = 171
    24 = 165
    25 = 98
    26 = 63
    27 = 48
    28 = 24
    29 = 53
    30 = 38
    31 = 28
    32 = 159
    33 = 60
    34 = 35
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 26).Value = $values[$row]
}
